$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

# Row 2
Set-TextValue "D2" "56.432.49"
$ws.Range("E2").Value = "  +2.27%  "

# Row 3
Set-TextValue "D3" "2.321.72"
$ws.Range("E3").Value = "  +1.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.38%  "

# Row 5
Set-TextValue "D5" "515.20"
$ws.Range("E5").Value = "  +1.84%  "

# Row 6
Set-TextValue "D6" "133.21"
$ws.Range("E6").Value = "  +3.13%  "

# Row 7
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  +0.40%  "

# Row 8
$ws.Range("E8").Value = "  +0.92%  "

# Row 9
Set-TextValue "D9" "2.343.63"
$ws.Range("E9").Value = "  +1.77%  "

# Row 10
$ws.Range("E10").Value = "  +5.95%  "

# Row 11
$ws.Range("E11").Value = "  +0.23%  "

# Row 12
$ws.Range("E12").Value = "  +5.77%  "

# Row 13
$ws.Range("E13").Value = "  -1.69%  "

# Row 14
Set-TextValue "D14" "2.760.26"
$ws.Range("E14").Value = "  +2.61%  "

# Row 15
Set-TextValue "D15" "23.60"
$ws.Range("E15").Value = "  +0.62%  "

# Row 16
Set-TextValue "D16" "56.608.11"
$ws.Range("E16").Value = "  +3.07%  "

# Row 17
$ws.Range("E17").Value = "  +1.72%  "

# Row 18
Set-TextValue "D18" "2.350.68"
$ws.Range("E18").Value = "  +3.03%  "

# Row 19
$ws.Range("E19").Value = "  +0.13%  "

# Row 20
$ws.Range("E20").Value = "  +2.02%  "

# Row 21
Set-TextValue "D21" "319.35"
$ws.Range("E21").Value = "  +3.79%  "

# Row 22
$ws.Range("E22").Value = "  +2.14%  "

# Row 23
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
Set-TextValue "D24" "60.64"
$ws.Range("E24").Value = "  +0.78%  "

# Row 25
Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  +0.50%  "

# Row 26
$ws.Range("E26").Value = "  +4.92%  "

# Row 27
$ws.Range("E27").Value = "  +3.20%  "

# Row 28
Set-TextValue "D28" "170.61"
$ws.Range("E28").Value = "  -0.34%  "

# Row 29
$ws.Range("E29").Value = "  +9.14%  "

# Row 30
$ws.Range("E30").Value = "  +4.04%  "

# Row 31
$ws.Range("E31").Value = "  +1.79%  "

# Row 32
$ws.Range("E32").Value = "  +2.19%  "

# Row 33
Set-TextValue "D33" "18.21"
$ws.Range("E33").Value = "  +1.28%  "

# Row 35
Set-TextValue "D35" "0.996"
$ws.Range("E35").Value = "  -0.12%  "

# Row 36
Set-TextValue "D36" "0.940"
$ws.Range("E36").Value = "  +2.26%  "

# Row 37
$ws.Range("E37").Value = "  +2.11%  "

# Row 38
$ws.Range("E38").Value = "  +3.92%  "

# Row 39
$ws.Range("E39").Value = "  +7.00%  "

# Row 40
Set-TextValue "D40" "37.43"
$ws.Range("E40").Value = "  +2.62%  "

# Row 42
Set-TextValue "D42" "137.60"
$ws.Range("E42").Value = "  +8.33%  "

# Row 43
$ws.Range("E43").Value = "  +3.92%  "

# Row 44
Set-TextValue "D44" "275.05"
$ws.Range("E44").Value = "  +9.75%  "

# Row 45
$ws.Range("E45").Value = "  -0.66%  "

# Row 46
Set-TextValue "D46" "0.0928"
$ws.Range("E46").Value = "  +2.76%  "

# Row 47
$ws.Range("E47").Value = "  +0.43%  "

# Row 48
Set-TextValue "D48" "0.558"
$ws.Range("E48").Value = "  +1.26%  "

# Row 49
$ws.Range("E49").Value = "  +3.96%  "

# Row 50
$ws.Range("E50").Value = "  +0.82%  "

# Row 51
Set-TextValue "D51" "16.70"
$ws.Range("E51").Value = "  +1.08%  "
